# Daten aktualisiert am 2024-03-03
#
# A new constituent (AMZN / Amazon) was added to the Dow Jones ticker list.
# It is inserted as row 5 (alphabetically, right after AMGN), which pushes
# every following company down by one row; the sheet still ends at row 31
# because the last old row (WBA, Walgreens Boots Alliance) drops off the
# bottom of the tracked range while WMT (row 31) stays put.
#
# Rather than using a real "insert row" (which would also shift row 31),
# we simply rewrite the values of rows 5-30 with the new, shifted data -
# this produces exactly the same end state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is brand new; rows 6-30 take over the data that used to sit one row
# above them (i.e. old row N -> new row N+1) for old rows 5..29.
$newData = @(
    @('AMZN', 'Amazon',                 'Retailing',                      '2024-02-26'),
    @('AAPL', 'Apple',                  'Information technology',         '2015-03-19'),
    @('BA',   'Boeing',                 'Aerospace and defense',          '1987-03-12'),
    @('CAT',  'Caterpillar',            'Construction and mining',        '1991-05-06'),
    @('CVX',  'Chevron',                'Petroleum industry',             '2008-02-19'),
    @('CSCO', 'Cisco',                  'Information technology',         '2009-06-08'),
    @('KO',   'Coca-Cola',              'Drink industry',                 '1987-03-12'),
    @('DIS',  'Disney',                 'Broadcasting and entertainment', '1991-05-06'),
    @('DOW',  'Dow',                    'Chemical industry',              '1991-05-06'),
    @('GS',   'Goldman Sachs',          'Financial services',             '2019-04-02'),
    @('HD',   'Home Depot',             'Home Improvement',               '1999-11-01'),
    @('HON',  'Honeywell',              'Conglomerate',                   '2020-08-31'),
    @('IBM',  'IBM',                    'Information technology',         '1979-06-29'),
    @('INTC', 'Intel',                  'Semiconductor industry',         '1999-11-01'),
    @('JNJ',  'Johnson & Johnson',      'Pharmaceutical industry',        '1997-03-17'),
    @('JPM',  'JPMorgan Chase',         'Financial services',             '1991-05-06'),
    @('MCD',  "McDonald's",             'Food industry',                  '1985-10-30'),
    @('MRK',  'Merck',                  'Pharmaceutical industry',        '1979-06-29'),
    @('MSFT', 'Microsoft',              'Information technology',         '1999-11-01'),
    @('NKE',  'Nike',                   'Clothing industry',              '2013-09-20'),
    @('PG',   'Procter & Gamble',       'Fast-moving consumer goods',     '1932-05-26'),
    @('CRM',  'Salesforce',             'Information technology',         '2020-08-31'),
    @('TRV',  'Travelers',              'Insurance',                      '2009-06-08'),
    @('UNH',  'UnitedHealth Group',     'Managed health care',            '2012-09-24'),
    @('VZ',   'Verizon',                'Telecommunications industry',    '2004-04-08'),
    @('V',    'Visa',                   'Financial services',             '2013-09-20')
)

$startRow = 5
$endRow = $startRow + $newData.Length - 1

# The "Date added" column holds plain text like "2024-02-26". Excel would
# normally auto-convert such a string into a date serial number, so the
# column is temporarily forced to Text format while the values are written.
$dateRangeAddress = "D" + $startRow + ":D" + $endRow
$dateRange = $ws.Range($dateRangeAddress)
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $values = $newData[$i]

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

# Put the number format / style back to the sheet default so the rewritten
# cells don't end up looking any different than the untouched ones.
$dateRange.NumberFormat = "General"
$dateRange.Style = "Normal"
